# Applies the diff: update tickers in column C (rows 2-5), clear the
# stray ticker in B2, and append three new data rows (6-8) with the same
# look/formatting as the existing rows (copied from row 5 so the "A"
# column keeps its bordered/centered/bold style and the other inlineStr
# cells keep the same empty-but-present shape).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing ticker cells (column C, rows 2-5) ---
$ws.Range("C2").Value = "NSE:ADL"
$ws.Range("C3").Value = "NSE:AMRUTANJAN"
$ws.Range("C4").Value = "NSE:AVALON"
$ws.Range("C5").Value = "NSE:BALKRISIND"

# --- B2 no longer carries a ticker ---
$ws.Range("B2").ClearContents()

# --- Append rows 6, 7, 8 ---
# Copy the full row-5 formatting/shape down so the new rows match the
# existing ones (bordered/centered/bold "A" style, empty placeholder
# cells for B/D/E/F), then overwrite the index + ticker values.
$ws.Range("A5:F5").Copy($ws.Range("A6:F6"))
$ws.Range("A6").Value = 4
$ws.Range("C6").Value = "NSE:DREAMFOLKS"

$ws.Range("A5:F5").Copy($ws.Range("A7:F7"))
$ws.Range("A7").Value = 5
$ws.Range("C7").Value = "NSE:JLHL"

$ws.Range("A5:F5").Copy($ws.Range("A8:F8"))
$ws.Range("A8").Value = 6
$ws.Range("C8").Value = "NSE:LGHL"
